# Reorder the "Recorded By" (column G) values so that any "System"/"system"
# token(s) come first, followed by the remaining (email) tokens, preserving
# their relative order. Only touches cells that actually contain a "System"
# token alongside other tokens; cells with a single value, or with no
# "System" token at all, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "*,*") { continue }
    if ($val -notmatch "(?i)(^|,\s*)system(\s*,|$)") { continue }

    $parts = $val -split ','
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $systemParts = @($trimmed | Where-Object { $_.ToLower() -eq 'system' })
    $otherParts = @($trimmed | Where-Object { $_.ToLower() -ne 'system' })

    $newParts = @($systemParts + $otherParts)
    $newVal = [string]::Join(', ', $newParts)

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
